$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 73: Person - Lockout No Access Users (mirrors row 69, Organization - Lockout No Access Users) ---
$ws.Range("B69:G69").Copy() | Out-Null
$ws.Range("B73:G73").PasteSpecial(-4122) | Out-Null
$ws.Range("C73").Value = "PERSON"
$ws.Range("B73").Value = "Person – Lockout No Access Users"
$ws.Range("G73").Value = "mandatory deny read to No Access"
$ws.Rows.Item(73).RowHeight = 30

# --- Row 74: Person - default read access (mirrors row 70, Organization - default read access) ---
$ws.Range("B70:G70").Copy() | Out-Null
$ws.Range("B74:G74").PasteSpecial(-4122) | Out-Null
$ws.Range("B74").Value = "Person – default read access"
$ws.Range("C74").Value = "PERSON"
$ws.Range("G74").Value = "grant read to owner, owning group, collaborator, reader, *"
$ws.Rows.Item(74).RowHeight = 30

# --- Row 75: Person - Only participants can save (mirrors row 71, Organization - Only participants can save) ---
$ws.Range("B71:G71").Copy() | Out-Null
$ws.Range("B75:G75").PasteSpecial(-4122) | Out-Null
$ws.Range("B75").Value = "Person – Only participants can save"
$ws.Range("C75").Value = "PERSON"
$ws.Range("G75").Value = "grant save to owner, owning group, collaborator"
$ws.Rows.Item(75).RowHeight = 30

# --- Row 76: Person - Restricted Flag (mirrors row 72, Organization - Restricted Flag, ---
# --- except column C keeps the un-bordered "17" style used on rows 73-75, not row 72's "18") ---
$ws.Range("B72:G72").Copy() | Out-Null
$ws.Range("B76:G76").PasteSpecial(-4122) | Out-Null
$ws.Range("C69").Copy() | Out-Null
$ws.Range("C76").PasteSpecial(-4122) | Out-Null
$ws.Range("B76").Value = "Person – Restricted Flag"
$ws.Range("C76").Value = "PERSON"
$ws.Range("D76").Value = "restricted"
$ws.Range("G76").Value = "deny read to *"

$excel.CutCopyMode = $false

# Update the active selection to match the new end of the table
$ws.Range("B77").Select() | Out-Null
